$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): marks per right answer and penalty per wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the "Total" row (row 12): total score and the "scored / max" summary text
$ws.Range("B12").Value = 84
$ws.Range("E12").Value = "84 / 112"
